$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B7: was a text date "01.04.2025"; becomes a real date serial (45754 = 2025-04-07)
# with a short-date number format (builtin numFmtId 14).
$ws.Range("B7").Value = 45754
$ws.Range("B7").NumberFormat = "mm-dd-yy"

# Selection moved to E12 (next active cell clicked by the user).
$null = $ws.Range("E12").Select()
